# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve rows
# in each crafting-job sheet. Values mirror the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 273.2619
$ws.Range("I33").Value = 163.94737
$ws.Range("J33").Value = 1311.75
$ws.Range("K33").Value = 163.94737
$ws.Range("L33").Value = 1311.75
$ws.Range("M33").Value = 65.05262999999999
$ws.Range("N33").Value = -1769.75

$ws.Range("H108").Value = 45280
$ws.Range("J108").Value = 45280
$ws.Range("L108").Value = 45280
$ws.Range("N108").Value = -52960

$ws.Range("H109").Value = 36529.75
$ws.Range("J109").Value = 36529.75
$ws.Range("L109").Value = 36529.75
$ws.Range("N109").Value = -39303.75

$ws.Range("H117").Value = 48727.8
$ws.Range("J117").Value = 48727.8
$ws.Range("L117").Value = 48727.8
$ws.Range("N117").Value = -57905.8

$ws.Range("H132").Value = 12781.169
$ws.Range("I132").Value = 1931.5286
$ws.Range("J132").Value = 71202.30499999999
$ws.Range("K132").Value = 5794.585800000001
$ws.Range("L132").Value = 213606.915
$ws.Range("M132").Value = -3264.585800000001
$ws.Range("N132").Value = -218666.915

$ws.Range("H137").Value = 3986.087
$ws.Range("I137").Value = 1230.6666
$ws.Range("J137").Value = 6992
$ws.Range("K137").Value = 3691.9998
$ws.Range("L137").Value = 20976
$ws.Range("M137").Value = -1141.9998
$ws.Range("N137").Value = -26076

$ws.Range("H138").Value = 1895.5955
$ws.Range("I138").Value = 1037.849
$ws.Range("J138").Value = 3158.389
$ws.Range("K138").Value = 3113.547
$ws.Range("L138").Value = 9475.167000000001
$ws.Range("M138").Value = 2026.453
$ws.Range("N138").Value = -19755.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 41410
$ws.Range("J109").Value = 41410
$ws.Range("L109").Value = 41410
$ws.Range("N109").Value = -44184

$ws.Range("H118").Value = 49401
$ws.Range("J118").Value = 49401
$ws.Range("L118").Value = 49401
$ws.Range("N118").Value = -52715

$ws.Range("H138").Value = 44500
$ws.Range("J138").Value = 44500
$ws.Range("L138").Value = 44500
$ws.Range("N138").Value = -54780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 48847.668
$ws.Range("J126").Value = 48847.668
$ws.Range("L126").Value = 48847.668
$ws.Range("N126").Value = -58727.668

$ws.Range("H132").Value = 33816
$ws.Range("J132").Value = 33816
$ws.Range("L132").Value = 33816
$ws.Range("N132").Value = -43936

$ws.Range("H133").Value = 45499.75
$ws.Range("J133").Value = 45499.75
$ws.Range("L133").Value = 45499.75
$ws.Range("N133").Value = -55619.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 16000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 16000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 16000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -16470

$ws.Range("H31").Value = 2800.9062
$ws.Range("I31").Value = 1185.8077
$ws.Range("J31").Value = 3400.8
$ws.Range("K31").Value = 1185.8077
$ws.Range("L31").Value = 3400.8
$ws.Range("M31").Value = -890.8077000000001
$ws.Range("N31").Value = -3990.8

$ws.Range("H34").Value = 2800.9062
$ws.Range("I34").Value = 1185.8077
$ws.Range("J34").Value = 3400.8
$ws.Range("K34").Value = 1185.8077
$ws.Range("L34").Value = 3400.8
$ws.Range("M34").Value = -983.8077000000001
$ws.Range("N34").Value = -3804.8

$ws.Range("H58").Value = 1281.942
$ws.Range("I58").Value = 1056.9348
$ws.Range("J58").Value = 1731.9565
$ws.Range("K58").Value = 1056.9348
$ws.Range("L58").Value = 1731.9565
$ws.Range("M58").Value = -853.9348
$ws.Range("N58").Value = -2137.9565

$ws.Range("H131").Value = 41993
$ws.Range("J131").Value = 41993
$ws.Range("L131").Value = 41993
$ws.Range("N131").Value = -52073

$ws.Range("H136").Value = 1281.942
$ws.Range("I136").Value = 1056.9348
$ws.Range("J136").Value = 1731.9565
$ws.Range("K136").Value = 3170.8044
$ws.Range("L136").Value = 5195.8695
$ws.Range("M136").Value = -620.8044
$ws.Range("N136").Value = -10295.8695

$ws.Range("H137").Value = 67685.57000000001
$ws.Range("J137").Value = 67685.57000000001
$ws.Range("L137").Value = 67685.57000000001
$ws.Range("N137").Value = -77885.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 143.75
$ws.Range("I12").Value = 156.33333
$ws.Range("J12").Value = 134.3125
$ws.Range("K12").Value = 468.99999
$ws.Range("L12").Value = 402.9375
$ws.Range("M12").Value = -295.99999
$ws.Range("N12").Value = -748.9375

$ws.Range("H122").Value = 8709.385
$ws.Range("J122").Value = 13810.625
$ws.Range("L122").Value = 124295.625
$ws.Range("N122").Value = -129195.625

$ws.Range("H131").Value = 10411.583
$ws.Range("J131").Value = 3071.4285
$ws.Range("L131").Value = 9214.2855
$ws.Range("N131").Value = -19294.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10270.308
$ws.Range("I126").Value = 13690.444
$ws.Range("J126").Value = 2575
$ws.Range("K126").Value = 41071.33199999999
$ws.Range("L126").Value = 7725
$ws.Range("M126").Value = -38601.33199999999
$ws.Range("N126").Value = -12665

$ws.Range("H135").Value = 35290
$ws.Range("J135").Value = 35290
$ws.Range("L135").Value = 35290
$ws.Range("N135").Value = -45430

$ws.Range("H138").Value = 48000
$ws.Range("J138").Value = 48000
$ws.Range("L138").Value = 48000
$ws.Range("N138").Value = -58280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1682.1082
$ws.Range("I16").Value = 1492.1515
$ws.Range("J16").Value = 3249.25
$ws.Range("K16").Value = 1492.1515
$ws.Range("L16").Value = 3249.25
$ws.Range("M16").Value = -1322.1515
$ws.Range("N16").Value = -3589.25

$ws.Range("H134").Value = 51666.332
$ws.Range("J134").Value = 51666.332
$ws.Range("L134").Value = 51666.332
$ws.Range("N134").Value = -61806.332

$ws.Range("H135").Value = 49333.332
$ws.Range("J135").Value = 49333.332
$ws.Range("L135").Value = 49333.332
$ws.Range("N135").Value = -59473.332

$ws.Range("H136").Value = 1248.4717
$ws.Range("I136").Value = 941.57776
$ws.Range("K136").Value = 2824.73328
$ws.Range("M136").Value = -274.7332799999999

$ws.Range("H137").Value = 39900
$ws.Range("J137").Value = 39900
$ws.Range("L137").Value = 39900
$ws.Range("N137").Value = -50100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 44694
$ws.Range("J119").Value = 44694
$ws.Range("L119").Value = 44694
$ws.Range("N119").Value = -54370

$ws.Range("H136").Value = 263042.84
$ws.Range("I136").Value = 524624.9
$ws.Range("J136").Value = 1460.8572
$ws.Range("K136").Value = 1573874.7
$ws.Range("L136").Value = 4382.571599999999
$ws.Range("M136").Value = -1571324.7
$ws.Range("N136").Value = -9482.571599999999

$ws.Range("H138").Value = 39087.777
$ws.Range("L138").Value = 39087.777
$ws.Range("N138").Value = -49367.777
